$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 1.43837993119209
$ws.Range("B3").Value = 1.3354509335489999
$ws.Range("B4").Value = 213381.57839400001
$ws.Range("A42").Value = 213381.57839400001

$ws.Range("B41").Value = 0.61896070999999997
$ws.Range("C41").Value = 0.74816643999999999
$ws.Range("D41").Value = 0.86089705999999999
$ws.Range("E41").Value = 0.97412798
$ws.Range("F41").Value = 1.0735994
$ws.Range("G41").Value = 1.4520608100000001
$ws.Range("H41").Value = 2.55346954

$wb.Application.Calculate()
